$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the marking/total marks values on the marksheet
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"

$wb.Save()
